$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("experiment_1")

# Replace curly quotes with straight quotes in the "Description" column
# (B3, "Create a room with..." has no curly quotes, so it is left untouched)
$ws.Range("B2").Value = 'Add a wall "Test Wall 1" and a "Test Door" in the building'
$ws.Range("B4").Value = 'Name the building "Residential Building"'
$ws.Range("B5").Value = 'Change the name of the TestWall, to "Base Wall" and to its door to "Base Door"'
$ws.Range("B6").Value = 'Remove the door from the "Base Wall"'
$ws.Range("B7").Value = 'Delete all walls that creates with "Base Wall" a room'

# Update the view: scroll position, selection, and zoom
$ws.Application.ActiveWindow.Zoom = 101
$ws.Range("B7").Select()
